$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    3 = @(0.1169995834814548, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.426980108624251)
    4 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    5 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455)
    6 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.178645819794754)
    7 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    8 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.536033448013082)
    9 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
